$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# Add the "NA" value under duplicate_image_filename (column E) for rows 2-21
$ws.Range("E2:E21").Value = "NA"

# Re-touch F1 (pre-existing blank placeholder cell) so the round trip keeps
# it truly empty instead of resurrecting a stray cached shared-string value.
$ws.Range("F1").Value = ""
